$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.19726943969727
$ws.Range("B3").Value = 7.934562683105469
$ws.Range("B4").Value = 5.874011516571045
$ws.Range("B5").Value = 6.663724422454834
$ws.Range("B6").Value = 7.561977863311768
$ws.Range("B7").Value = 11.536208152771
$ws.Range("B8").Value = 13.69966888427734
$ws.Range("B9").Value = 10.80807781219482
$ws.Range("B10").Value = 4.527147769927979
$ws.Range("B11").Value = 5.923291206359863
$ws.Range("B12").Value = 4.451257228851318
$ws.Range("B13").Value = 7.551928043365479
$ws.Range("B14").Value = 10.58013725280762
$ws.Range("B15").Value = 10.74296760559082
$ws.Range("B16").Value = 10.26916980743408
$ws.Range("B17").Value = 6.098263263702393
$ws.Range("B18").Value = 4.29344367980957
$ws.Range("B19").Value = 13.55679416656494
$ws.Range("B20").Value = 19.59993553161621
$ws.Range("B21").Value = 6.560589790344238
$ws.Range("B22").Value = 5.0816330909729
$ws.Range("B23").Value = 4.274739742279053
$ws.Range("B24").Value = 6.05767822265625
$ws.Range("B25").Value = 6.501260280609131
$ws.Range("B26").Value = 7.32531213760376
$ws.Range("B27").Value = 12.09658622741699
$ws.Range("B28").Value = 8.651673316955566
$ws.Range("B29").Value = 11.36209487915039
$ws.Range("B30").Value = 14.84198379516602
$ws.Range("B31").Value = 11.55424976348877
$ws.Range("B32").Value = 16.05239677429199
$ws.Range("B33").Value = 6.340110778808594
$ws.Range("B34").Value = 19.92737579345703
$ws.Range("B35").Value = 25.71750259399414
$ws.Range("B36").Value = 14.04750919342041
$ws.Range("B37").Value = 7.858399868011475
$ws.Range("B38").Value = 5.839068412780762
$ws.Range("B39").Value = 6.640105247497559
$ws.Range("B40").Value = 7.540596961975098
$ws.Range("B41").Value = 11.44162750244141
$ws.Range("B42").Value = 13.59792041778564
$ws.Range("B43").Value = 10.69638729095459
$ws.Range("B44").Value = 4.536865234375
$ws.Range("B45").Value = 5.864860057830811
$ws.Range("B46").Value = 4.458763122558594
$ws.Range("B47").Value = 7.484989643096924
$ws.Range("B48").Value = 10.48827266693115
$ws.Range("B49").Value = 10.67738342285156
$ws.Range("B50").Value = 10.18572998046875
$ws.Range("B51").Value = 6.018743515014648
$ws.Range("B52").Value = 4.29326868057251
$ws.Range("B53").Value = 13.44388008117676
$ws.Range("B54").Value = 19.38541984558105
$ws.Range("B55").Value = 6.491823196411133
$ws.Range("B56").Value = 5.036499500274658
$ws.Range("B57").Value = 4.248396873474121
$ws.Range("B58").Value = 5.995769500732422
$ws.Range("B59").Value = 6.410431385040283
$ws.Range("B60").Value = 7.259796142578125
$ws.Range("B61").Value = 11.92684364318848
$ws.Range("B62").Value = 8.585808753967285
$ws.Range("B63").Value = 11.28206729888916
$ws.Range("B64").Value = 14.72089862823486
$ws.Range("B65").Value = 11.45261573791504
$ws.Range("B66").Value = 15.84002017974854
$ws.Range("B67").Value = 6.196976661682129
$ws.Range("B68").Value = 19.68788909912109
$ws.Range("B69").Value = 25.40290832519531
$ws.Range("B70").Value = 13.89502143859863
$ws.Range("B71").Value = 7.770692825317383
$ws.Range("B72").Value = 5.808511734008789
$ws.Range("B73").Value = 6.605836868286133
$ws.Range("B74").Value = 7.507067203521729
$ws.Range("B75").Value = 11.33576107025146
$ws.Range("B76").Value = 13.48910808563232
$ws.Range("B77").Value = 10.56950664520264
$ws.Range("B78").Value = 4.533520221710205
$ws.Range("B79").Value = 5.803475379943848
$ws.Range("B80").Value = 4.458601474761963
$ws.Range("B81").Value = 7.420427799224854
$ws.Range("B82").Value = 10.39489078521729
$ws.Range("B83").Value = 10.60312271118164
$ws.Range("B84").Value = 10.1016731262207
$ws.Range("B85").Value = 5.935788631439209
$ws.Range("B86").Value = 4.278801918029785
$ws.Range("B87").Value = 13.31656265258789
$ws.Range("B88").Value = 19.15806579589844
$ws.Range("B89").Value = 6.417155742645264
$ws.Range("B90").Value = 4.981106281280518
$ws.Range("B91").Value = 4.218957424163818
$ws.Range("B92").Value = 5.928547859191895
$ws.Range("B93").Value = 6.316637992858887
$ws.Range("B94").Value = 7.186127662658691
$ws.Range("B95").Value = 11.76066780090332
$ws.Range("B96").Value = 8.511568069458008
$ws.Range("B97").Value = 11.19528484344482
$ws.Range("B98").Value = 14.59385108947754
$ws.Range("B99").Value = 11.34084892272949
$ws.Range("B100").Value = 15.6251916885376
$ws.Range("B101").Value = 6.07497501373291
$ws.Range("B102").Value = 19.44630241394043
$ws.Range("B103").Value = 25.08788681030273
$ws.Range("B104").Value = 13.73633766174316
$ws.Range("B105").Value = 7.671136379241943
$ws.Range("B106").Value = 5.777669906616211
$ws.Range("B107").Value = 6.559080600738525
$ws.Range("B108").Value = 7.459660530090332
$ws.Range("B109").Value = 11.21639442443848
$ws.Range("B110").Value = 13.36867618560791
$ws.Range("B111").Value = 10.42665672302246
$ws.Range("B112").Value = 4.514937877655029
$ws.Range("B113").Value = 5.735990047454834
$ws.Range("B114").Value = 4.447064876556396
$ws.Range("B115").Value = 7.355311393737793
$ws.Range("B116").Value = 10.29762840270996
$ws.Range("B117").Value = 10.52012538909912
$ws.Range("B118").Value = 10.01451015472412
$ws.Range("B119").Value = 5.851010322570801
$ws.Range("B120").Value = 4.250696659088135
$ws.Range("B121").Value = 13.18251800537109
$ws.Range("B122").Value = 18.92146301269531
$ws.Range("B123").Value = 6.337779521942139
$ws.Range("B124").Value = 4.916608333587646
$ws.Range("B125").Value = 4.189223289489746
$ws.Range("B126").Value = 5.859615325927734
$ws.Range("B127").Value = 6.228358268737793
$ws.Range("B128").Value = 7.105429172515869
$ws.Range("B129").Value = 11.59707832336426
$ws.Range("B130").Value = 8.436002731323242
$ws.Range("B131").Value = 11.10911083221436
$ws.Range("B132").Value = 14.46375751495361
$ws.Range("B133").Value = 11.2252779006958
$ws.Range("B134").Value = 15.41115760803223
$ws.Range("B135").Value = 5.967560768127441
$ws.Range("B136").Value = 19.20959281921387
$ws.Range("B137").Value = 24.77530097961426
